$d = $word.ActiveDocument

# 1) "Objetivos" paragraph: insert a manual line break between
#    "...competências no" and "desenvolvimento de projetos..."
$old1 = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências nodesenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de metais e a redução de ocorrência de falhas estruturais baseado no trinômio propriedades, estrutura metalúrgica e processamento metalúrgico dos metais aplicado a engenharia permitindo aos alunos a prática da redação científica e da busca de projetos para incentivar a solução de problemas em engenharia."
$new1 = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências no^ldesenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de metais e a redução de ocorrência de falhas estruturais baseado no trinômio propriedades, estrutura metalúrgica e processamento metalúrgico dos metais aplicado a engenharia permitindo aos alunos a prática da redação científica e da busca de projetos para incentivar a solução de problemas em engenharia."
$r1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) "Programa" body paragraph: insert a manual line break between
#    "Programa" and "1.INTRODUÇÃO AO CONCEITO..."
$old2 = "Programa1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais metálicos. Comportamento elástico e plástico de metais suas ligas e materiais não ferrosos. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação. Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Estudo de ligas não metálicas. Tratamentos térmicos em aços e ligas especiais. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência. Impacto e a transição dúctil-frágil. 6. Influência da temperatura sobre o comportamento mecânico dos metais. Aspectos básicos da análise de falhas em materiais metálicos."
$new2 = "Programa^l1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais metálicos. Comportamento elástico e plástico de metais suas ligas e materiais não ferrosos. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação. Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Estudo de ligas não metálicas. Tratamentos térmicos em aços e ligas especiais. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência. Impacto e a transição dúctil-frágil. 6. Influência da temperatura sobre o comportamento mecânico dos metais. Aspectos básicos da análise de falhas em materiais metálicos."
$r2 = $d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3) "Bibliografia" paragraph: insert a manual line break before each
#    numbered reference (2. through 11.)
$old3 = "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009.2. Hearn, E.J. Mechanics of Materials: An Introduction to the Mechanics of Elastic and Plastic Deformation of Solids and Structural Components, Pergamon Press, 1985.3. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1981.4. Hull, D. Introduction to Dislocations, Pergamon Press, 1965.5. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967.6. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.7. Van Vlack, L.H. Princípios de Ciência dos materiais, Ed. Edgard Blucher Ltda., 1970.8. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008.9. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall,1988.10. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008.11. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
$new3 = "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009.^l2. Hearn, E.J. Mechanics of Materials: An Introduction to the Mechanics of Elastic and Plastic Deformation of Solids and Structural Components, Pergamon Press, 1985.^l3. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1981.^l4. Hull, D. Introduction to Dislocations, Pergamon Press, 1965.^l5. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967.^l6. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.^l7. Van Vlack, L.H. Princípios de Ciência dos materiais, Ed. Edgard Blucher Ltda., 1970.^l8. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008.^l9. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall,1988.^l10. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008.^l11. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
$r3 = $d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

Write-Output "Objetivos replace: $r1"
Write-Output "Programa replace: $r2"
Write-Output "Bibliografia replace: $r3"
